$d = $word.ActiveDocument

# --- Edit 1: Overview body paragraph - split out "OpenSCAD" into its own run
# wrapped in proofErr spellStart/spellEnd markers (simulating Word's spell-check
# flagging of the non-dictionary word), same as done for Build Instructions body.
$find1 = $d.Content.Find
$find1.ClearFormatting()
$find1.Text = "The Cup Holder Resizing Ring is intended"
$found1 = $find1.Execute()
if (-not $found1) { throw "Could not locate Overview paragraph" }
$p1 = $find1.Parent.Paragraphs(1).Range

$xml1 = '<?xml version="1.0"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml"><w:body><w:p w14:paraId="02EB378F" w14:textId="5E1B2743" w:rsidR="00855181" w:rsidRDefault="39BEF063" w:rsidP="0E3361DC"><w:pPr><w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr><w:r w:rsidRPr="0E3361DC"><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve">The Cup Holder Resizing Ring is intended to reduce the diameter of an existing cup holder to better fit cups or mugs and prevent them from tipping over and falling out. The design is made in </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>OpenSCAD</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve"> and is fully parametric so that with 4 to 5 measurements, users can receive a custom sized ring perfect for their specific needs. The Resizing Ring is as simple to use as dropping into the cup holder and then placing your drink inside.</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$p1.InsertXML($xml1)

# --- Edit 2: Build Instructions body paragraph - same OpenSCAD split.
$find2 = $d.Content.Find
$find2.ClearFormatting()
$find2.Text = "Take the desired dimensions"
$found2 = $find2.Execute()
if (-not $found2) { throw "Could not locate Build Instructions paragraph" }
$p2 = $find2.Parent.Paragraphs(1).Range

$xml2 = '<?xml version="1.0"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml"><w:body><w:p w14:paraId="1D4D1203" w14:textId="2ECE4BAD" w:rsidR="00D744A5" w:rsidRPr="00D744A5" w:rsidRDefault="00D744A5" w:rsidP="00D744A5"><w:pPr><w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr><w:r w:rsidRPr="0E3361DC"><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve">Take the desired dimensions and input them into the </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>OpenSCAD</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve"> file to generate the STL. Then print the file and the Resizing Ring is ready to use.</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$p2.InsertXML($xml2)

Write-Output "Edits applied"
